$d = $word.ActiveDocument

# The second paragraph currently ends with "...they are not equal with the others."
# We need to append three more runs of text to that same paragraph:
#   1) "Actually,...outstand.Don"   (rFonts hint="eastAsia")
#   2) "'"                           (rFonts, no hint - matches existing curly-quote runs)
#   3) "t be jealous other blindly." (rFonts hint="eastAsia")

$p = $d.Paragraphs.Item(2)
$r = $p.Range

# Source ranges already present in the document whose run formatting we want to
# replicate exactly (including the w:hint="eastAsia" attribute pattern used
# throughout this paragraph).
$hintSrc = $d.Range(12, 13)    # "I" from "I believe..." -> rFonts w:hint="eastAsia"
$quoteSrc = $d.Range(104, 105) # "'" from "It's" -> rFonts without hint, the curly apostrophe run

$part1 = "Actually,everyone has something special,there is no need to envy others,we might be envied by others.We can just be ourselves,just do what we want to do.Today,more and more people are doing plastic surgery,because they are not confident with themselves,they think the beautiful face makes them outstand,this is wrong idea.Just be yourself,you have the shinning point that makes you outstand.Don"
$part2 = [string][char]0x2019
$part3 = "t be jealous other blindly."

# --- Run 1 ---
$insertStart = $r.End - 1
$r.InsertAfter("X")
$newRange = $d.Range($insertStart, $insertStart + 1)
$newRange.FormattedText = $hintSrc.FormattedText
$newRange2 = $d.Range($insertStart, $insertStart + 1)
$newRange2.Text = $part1

# --- Run 2 ---
$r = $p.Range
$insertStart = $r.End - 1
$r.InsertAfter("X")
$newRange = $d.Range($insertStart, $insertStart + 1)
$newRange.FormattedText = $quoteSrc.FormattedText
$newRange2 = $d.Range($insertStart, $insertStart + 1)
$newRange2.Text = $part2

# --- Run 3 ---
$r = $p.Range
$insertStart = $r.End - 1
$r.InsertAfter("X")
$newRange = $d.Range($insertStart, $insertStart + 1)
$newRange.FormattedText = $hintSrc.FormattedText
$newRange2 = $d.Range($insertStart, $insertStart + 1)
$newRange2.Text = $part3

Write-Output $d.Content.Text
